# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) timestamps for the zh-cn and de-de
# report sheets to reflect the latest report-generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 09:08:15"
$wsZhCn.Range("H2").Value = "2016-03-13 09:08:33"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 09:08:19"
$wsDeDe.Range("H2").Value = "2016-03-13 09:08:39"
